# Script: applies the "Atualizado por script em 11-11-2023 08:45" update
# to the spain_laliga_2023-2024 sheet.
#
# Summary of the change:
#   1. Rows 78 and 79 had their match details (columns F:V) swapped back
#      to the correct order (Betis-Valencia vs Atl.Madrid-Cadiz).
#   2. Rows 88 and 89 had their match details (columns F:V) swapped back
#      to the correct order (Alaves-Betis vs Celta Vigo-Getafe).
#   3. A new match (Ath Bilbao 4-3 Celta Vigo) was appended as row 121.
#
# Values are written as literals (rather than read-then-swapped) to avoid
# any floating-point round-trip drift through the COM Value/Text bridge.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 78: Betis 3-0 Valencia --------------------------------------
$ws.Cells.Item(78, 6).Value = "Betis"
$ws.Cells.Item(78, 8).Value = "Valencia"
$ws.Cells.Item(78, 9).Value = 0
$ws.Cells.Item(78, 10).Value = 2.5
$ws.Cells.Item(78, 11).Value = "24/09/2023 17:02"
$ws.Cells.Item(78, 12).Value = 2.07
$ws.Cells.Item(78, 13).Value = "01/10/2023 20:54"
$ws.Cells.Item(78, 14).Value = 3.19
$ws.Cells.Item(78, 15).Value = "24/09/2023 17:02"
$ws.Cells.Item(78, 16).Value = 3.41
$ws.Cells.Item(78, 17).Value = "01/10/2023 20:54"
$ws.Cells.Item(78, 18).Value = 3.11
$ws.Cells.Item(78, 19).Value = "24/09/2023 17:02"
$ws.Cells.Item(78, 20).Value = 4.04
$ws.Cells.Item(78, 21).Value = "01/10/2023 20:57"
$ws.Cells.Item(78, 22).Value = "https://www.betexplorer.com/football/spain/laliga/betis-valencia/vukArZ2c/"

# --- Row 79: Atl. Madrid 3-2 Cadiz CF ----------------------------------
$ws.Cells.Item(79, 6).Value = "Atl. Madrid"
$ws.Cells.Item(79, 8).Value = "Cadiz CF"
$ws.Cells.Item(79, 9).Value = 2
$ws.Cells.Item(79, 10).Value = 1.35
$ws.Cells.Item(79, 11).Value = "21/09/2023 22:03"
$ws.Cells.Item(79, 12).Value = 1.34
$ws.Cells.Item(79, 13).Value = "01/10/2023 20:50"
$ws.Cells.Item(79, 14).Value = 5.24
$ws.Cells.Item(79, 15).Value = "21/09/2023 22:03"
$ws.Cells.Item(79, 16).Value = 5.31
$ws.Cells.Item(79, 17).Value = "01/10/2023 20:59"
$ws.Cells.Item(79, 18).Value = 9.25
$ws.Cells.Item(79, 19).Value = "21/09/2023 22:03"
$ws.Cells.Item(79, 20).Value = 10.48
$ws.Cells.Item(79, 21).Value = "01/10/2023 20:59"
$ws.Cells.Item(79, 22).Value = "https://www.betexplorer.com/football/spain/laliga/atl-madrid-cadiz/E1cOKVAj/"

# --- Row 88: Alaves 1-1 Betis ------------------------------------------
$ws.Cells.Item(88, 6).Value = "Alaves"
$ws.Cells.Item(88, 7).Value = 1
$ws.Cells.Item(88, 8).Value = "Betis"
$ws.Cells.Item(88, 9).Value = 1
$ws.Cells.Item(88, 10).Value = 2.8
$ws.Cells.Item(88, 11).Value = "01/10/2023 20:24"
$ws.Cells.Item(88, 12).Value = 2.57
$ws.Cells.Item(88, 13).Value = "08/10/2023 18:28"
$ws.Cells.Item(88, 14).Value = 3.04
$ws.Cells.Item(88, 15).Value = "01/10/2023 20:24"
$ws.Cells.Item(88, 16).Value = 3.25
$ws.Cells.Item(88, 17).Value = "08/10/2023 18:28"
$ws.Cells.Item(88, 18).Value = 2.72
$ws.Cells.Item(88, 19).Value = "01/10/2023 20:24"
$ws.Cells.Item(88, 20).Value = 3.06
$ws.Cells.Item(88, 21).Value = "08/10/2023 18:22"
$ws.Cells.Item(88, 22).Value = "https://www.betexplorer.com/football/spain/laliga/alaves-betis/YNPlfW19/"

# --- Row 89: Celta Vigo 2-2 Getafe --------------------------------------
$ws.Cells.Item(89, 6).Value = "Celta Vigo"
$ws.Cells.Item(89, 7).Value = 2
$ws.Cells.Item(89, 8).Value = "Getafe"
$ws.Cells.Item(89, 9).Value = 2
$ws.Cells.Item(89, 10).Value = 1.92
$ws.Cells.Item(89, 11).Value = "28/09/2023 15:02"
$ws.Cells.Item(89, 12).Value = 2.04
$ws.Cells.Item(89, 13).Value = "08/10/2023 18:29"
$ws.Cells.Item(89, 14).Value = 3.25
$ws.Cells.Item(89, 15).Value = "28/09/2023 15:02"
$ws.Cells.Item(89, 16).Value = 3.29
$ws.Cells.Item(89, 17).Value = "08/10/2023 18:27"
$ws.Cells.Item(89, 18).Value = 4.85
$ws.Cells.Item(89, 19).Value = "28/09/2023 15:02"
$ws.Cells.Item(89, 20).Value = 4.38
$ws.Cells.Item(89, 21).Value = "08/10/2023 18:29"
$ws.Cells.Item(89, 22).Value = "https://www.betexplorer.com/football/spain/laliga/celta-vigo-getafe/0ARtdhXd/"

# --- New row 121: Ath Bilbao 4-3 Celta Vigo -----------------------------
# Copy formats from the last existing row (120) first, so the new row's
# index/date-time columns (A, E) keep the expected number formats/styles.
$ws.Range("A120:V120").Copy()
$ws.Range("A121:V121").PasteSpecial(-4122)

$ws.Cells.Item(121, 1).Value = 120
$ws.Cells.Item(121, 2).Value = "spain"
$ws.Cells.Item(121, 3).Value = "laliga"
$ws.Cells.Item(121, 4).Value = "2023-2024"
$ws.Cells.Item(121, 5).Value = 45240.875
$ws.Cells.Item(121, 6).Value = "Ath Bilbao"
$ws.Cells.Item(121, 7).Value = 4
$ws.Cells.Item(121, 8).Value = "Celta Vigo"
$ws.Cells.Item(121, 9).Value = 3
$ws.Cells.Item(121, 10).Value = 1.54
$ws.Cells.Item(121, 11).Value = "29/10/2023 11:02"
$ws.Cells.Item(121, 12).Value = 1.71
$ws.Cells.Item(121, 13).Value = "10/11/2023 20:57"
$ws.Cells.Item(121, 14).Value = 4.02
$ws.Cells.Item(121, 15).Value = "29/10/2023 11:02"
$ws.Cells.Item(121, 16).Value = 4.04
$ws.Cells.Item(121, 17).Value = "10/11/2023 20:59"
$ws.Cells.Item(121, 18).Value = 5.82
$ws.Cells.Item(121, 19).Value = "29/10/2023 11:02"
$ws.Cells.Item(121, 20).Value = 5.04
$ws.Cells.Item(121, 21).Value = "10/11/2023 20:59"
$ws.Cells.Item(121, 22).Value = "https://www.betexplorer.com/football/spain/laliga/ath-bilbao-celta-vigo/UPVu9gDU/"
